$d = $word.ActiveDocument

# The "www.plainsbororotary.org" hyperlink lives in the document's footer.
# Walk every section's footers (and headers, just in case) looking for it,
# turn the live hyperlink field back into plain (non-linked) text, and
# restyle the run to the new "not-a-link" look: navy (002060) text,
# no underline - matching the rest of the footer's contact-info run.
$targetColor = 6299648   # RGB(0,32,96) == hex 002060, Word's BGR-ish long value
$targetText = "www.plainsbororotary.org"

foreach ($sec in $d.Sections) {

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            while ($ftr.Range.Hyperlinks.Count -gt 0) {
                $ftr.Range.Hyperlinks.Item(1).Delete()
            }

            $rng = $ftr.Range.Duplicate
            $found = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
            if ($found) {
                $rng.Font.Color = $targetColor
                $rng.Font.Underline = 0
            }
        }
    }

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            while ($hdr.Range.Hyperlinks.Count -gt 0) {
                $hdr.Range.Hyperlinks.Item(1).Delete()
            }

            $rng = $hdr.Range.Duplicate
            $found = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
            if ($found) {
                $rng.Font.Color = $targetColor
                $rng.Font.Underline = 0
            }
        }
    }
}
